# Scheduled market-data refresh: recompute currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on the affected Leve rows of each job sheet, per the latest price pull.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 511.35294
$ws.Range("I15").Value = 511.35294
$ws.Range("K15").Value = 1534.05882
$ws.Range("M15").Value = -1365.05882
# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 4094.7
$ws.Range("J40").Value = 6150
$ws.Range("L40").Value = 6150
$ws.Range("N40").Value = -6500
# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 5187.5
$ws.Range("I43").Value = 4000.3333
$ws.Range("J43").Value = 5899.8
$ws.Range("K43").Value = 4000.3333
$ws.Range("L43").Value = 5899.8
$ws.Range("M43").Value = -3931.3333
$ws.Range("N43").Value = -6037.8
# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 10477.261
$ws.Range("I64").Value = 6887.5557
$ws.Range("J64").Value = 12784.929
$ws.Range("K64").Value = 6887.5557
$ws.Range("L64").Value = 12784.929
$ws.Range("M64").Value = -6639.5557
$ws.Range("N64").Value = -13280.929
# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 10477.261
$ws.Range("I67").Value = 6887.5557
$ws.Range("J67").Value = 12784.929
$ws.Range("K67").Value = 6887.5557
$ws.Range("L67").Value = 12784.929
$ws.Range("M67").Value = -6029.5557
$ws.Range("N67").Value = -14500.929
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2242
$ws.Range("I137").Value = 2242
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6726
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -4176

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3: Skillet Labor / Bronze Skillet
$ws.Range("H3").Value = 839.1429000000001
$ws.Range("I3").Value = 515
$ws.Range("J3").Value = 1649.5
$ws.Range("K3").Value = 515
$ws.Range("L3").Value = 1649.5
$ws.Range("M3").Value = -400
$ws.Range("N3").Value = -1879.5
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 5605.967
$ws.Range("I32").Value = 4027.3333
$ws.Range("K32").Value = 4027.3333
$ws.Range("M32").Value = -3740.3333
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2185.111
$ws.Range("I132").Value = 2185.111
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6555.333
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -4025.333

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 678.125
$ws.Range("I22").Value = 689.2857
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 689.2857
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -516.2857
$ws.Range("N22").Value = -946
# Row 36: I Saw What You Did There / Iron Chocobotail Saw
$ws.Range("H36").Value = 2133.3333
$ws.Range("I36").Value = 2133.3333
$ws.Range("K36").Value = 2133.3333
$ws.Range("M36").Value = -1599.3333
# Row 76: Keep Up with the Mechanics / Titanium-barreled Arquebus
$ws.Range("H76").Value = 14833
$ws.Range("J76").Value = 14833
$ws.Range("L76").Value = 14833
$ws.Range("N76").Value = -15463
# Row 79: Unconventional Weaponry (L) / Titanium-barreled Arquebus
$ws.Range("H79").Value = 14833
$ws.Range("J79").Value = 14833
$ws.Range("L79").Value = 14833
$ws.Range("N79").Value = -17017
# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 10583.333
$ws.Range("I82").Value = 10583.333
$ws.Range("K82").Value = 10583.333
$ws.Range("M82").Value = -10200.333
# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 10583.333
$ws.Range("I85").Value = 10583.333
$ws.Range("K85").Value = 10583.333
$ws.Range("M85").Value = -9257.333000000001
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3627.889
$ws.Range("I86").Value = 4540.2
$ws.Range("J86").Value = 2487.5
$ws.Range("K86").Value = 4540.2
$ws.Range("L86").Value = 2487.5
$ws.Range("M86").Value = -3417.2
$ws.Range("N86").Value = -4733.5
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3627.889
$ws.Range("I89").Value = 4540.2
$ws.Range("J89").Value = 2487.5
$ws.Range("K89").Value = 22701
$ws.Range("L89").Value = 12437.5
$ws.Range("M89").Value = -17085
$ws.Range("N89").Value = -23669.5
# Row 106: Fire for Hire / Molybdenum Rimfire
$ws.Range("H106").Value = 31666.666
$ws.Range("J106").Value = 31666.666
$ws.Range("L106").Value = 31666.666
$ws.Range("N106").Value = -34190.666
# Row 133: Paring Is Caring / Mountain Chromite Hatchet
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 785.7
$ws.Range("I16").Value = 794.125
$ws.Range("J16").Value = 752
$ws.Range("K16").Value = 794.125
$ws.Range("L16").Value = 752
$ws.Range("M16").Value = -507.125
$ws.Range("N16").Value = -1326
# Row 41: The Lone Bowman / Oak Longbow
$ws.Range("H41").Value = 3874.5
$ws.Range("I41").Value = 3874.5
$ws.Range("K41").Value = 3874.5
$ws.Range("M41").Value = -3446.5
# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 11621161
$ws.Range("I86").Value = 23237656
$ws.Range("J86").Value = 4666.6665
$ws.Range("K86").Value = 23237656
$ws.Range("L86").Value = 4666.6665
$ws.Range("M86").Value = -23236533
$ws.Range("N86").Value = -6912.6665
# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 11621161
$ws.Range("I89").Value = 23237656
$ws.Range("J89").Value = 4666.6665
$ws.Range("K89").Value = 116188280
$ws.Range("L89").Value = 23333.3325
$ws.Range("M89").Value = -116182664
$ws.Range("N89").Value = -34565.3325
# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 7019.727
$ws.Range("I99").Value = 6149.75
$ws.Range("J99").Value = 9339.666999999999
$ws.Range("K99").Value = 6149.75
$ws.Range("L99").Value = 9339.666999999999
$ws.Range("M99").Value = -4651.75
$ws.Range("N99").Value = -12335.667
# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 785.7
$ws.Range("I113").Value = 794.125
$ws.Range("J113").Value = 752
$ws.Range("K113").Value = 794.125
$ws.Range("L113").Value = 752
$ws.Range("M113").Value = 1375.875
$ws.Range("N113").Value = -5092
# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 7019.727
$ws.Range("I126").Value = 6149.75
$ws.Range("J126").Value = 9339.666999999999
$ws.Range("K126").Value = 18449.25
$ws.Range("L126").Value = 28019.001
$ws.Range("M126").Value = -15979.25
$ws.Range("N126").Value = -32959.001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 58: The Big Red / Red Coral Necklace
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1097.8572
$ws.Range("I102").Value = 1030.8334
$ws.Range("K102").Value = 1030.8334
$ws.Range("M102").Value = 591.1666
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 7578.4
$ws.Range("I126").Value = 7578.4
$ws.Range("K126").Value = 22735.2
$ws.Range("M126").Value = -20265.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 30: Packing a Punch / Goatskin Cesti
$ws.Range("H30").Value = 1965.9166
$ws.Range("I30").Value = 326.54544
$ws.Range("K30").Value = 326.54544
$ws.Range("M30").Value = -218.54544
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 3949
$ws.Range("I40").Value = 3949
$ws.Range("K40").Value = 3949
$ws.Range("M40").Value = -3813
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 3036.3635
$ws.Range("I46").Value = 2950
$ws.Range("J46").Value = 3055.5557
$ws.Range("K46").Value = 2950
$ws.Range("L46").Value = 3055.5557
$ws.Range("M46").Value = -2762
$ws.Range("N46").Value = -3431.5557
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1934.091
$ws.Range("I93").Value = 1840.2858
$ws.Range("J93").Value = 2098.25
$ws.Range("K93").Value = 1840.2858
$ws.Range("L93").Value = 2098.25
$ws.Range("M93").Value = -592.2858000000001
$ws.Range("N93").Value = -4594.25
# Row 105: Thick and Thin / Gazelleskin Corselet of Scouting
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 7408.636
$ws.Range("I122").Value = 5679
$ws.Range("K122").Value = 17037
$ws.Range("M122").Value = -14587

